$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column D first (Newsweek / South Cheyenne Solar), then column B
# (texaselectricnews / Arroyo Solar Energy Storage Hybrid) so the AEUG
# column (old C) shifts left into B and the remaining layout matches
# the target (A: project refs, B: AEUG Union Solar info).
$ws.Range("D:D").Delete() | Out-Null
$ws.Range("B:B").Delete() | Out-Null

$ws.Range("A3").Select() | Out-Null
